$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2899
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H6").Value = 50
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 150
$ws.Range("M6").Value = -38
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H18").Value = 166666820
$ws.Range("I18").Value = 166666820
$ws.Range("K18").Value = 166666820
$ws.Range("M18").Value = -166666536
$ws.Range("H21").Value = 290200
$ws.Range("J21").Value = 290200
$ws.Range("L21").Value = 290200
$ws.Range("N21").Value = -291136
$ws.Range("H23").Value = 290200
$ws.Range("J23").Value = 290200
$ws.Range("L23").Value = 290200
$ws.Range("N23").Value = -290668
$ws.Range("H38").Value = 5666.6665
$ws.Range("I38").Value = 2833.3333
$ws.Range("J38").Value = 11333.333
$ws.Range("K38").Value = 8499.999899999999
$ws.Range("L38").Value = 33999.999
$ws.Range("M38").Value = -8127.999899999999
$ws.Range("N38").Value = -34743.999
$ws.Range("H51").Value = 6855.5713
$ws.Range("J51").Value = 9000
$ws.Range("L51").Value = 9000
$ws.Range("N51").Value = -9968
$ws.Range("H58").Value = 1903.4286
$ws.Range("I58").Value = 85
$ws.Range("J58").Value = 6449.5
$ws.Range("K58").Value = 255
$ws.Range("L58").Value = 19348.5
$ws.Range("M58").Value = -105
$ws.Range("N58").Value = -19648.5
$ws.Range("H106").Value = 33339408
$ws.Range("I106").Value = 47623584
$ws.Range("K106").Value = 47623584
$ws.Range("M106").Value = -47622953
$ws.Range("H138").Value = 8255.204
$ws.Range("J138").Value = 8988.604499999999
$ws.Range("L138").Value = 26965.8135
$ws.Range("N138").Value = -37245.8135
$ws.Range("H139").Value = 96999
$ws.Range("J139").Value = 96999
$ws.Range("L139").Value = 96999
$ws.Range("N139").Value = -107279
$ws.Range("H140").Value = 60627.777
$ws.Range("J140").Value = 59456.25
$ws.Range("L140").Value = 59456.25
$ws.Range("N140").Value = -69816.25

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2662.5
$ws.Range("I32").Value = 2165.6667
$ws.Range("K32").Value = 2165.6667
$ws.Range("M32").Value = -1878.6667
$ws.Range("H43").Value = 9336.666999999999
$ws.Range("J43").Value = 9336.666999999999
$ws.Range("L43").Value = 9336.666999999999
$ws.Range("N43").Value = -9962.666999999999
$ws.Range("H45").Value = 2433.2856
$ws.Range("I45").Value = 1280.1333
$ws.Range("K45").Value = 1280.1333
$ws.Range("M45").Value = -903.1333
$ws.Range("H122").Value = 695499.75
$ws.Range("I122").Value = 1836666
$ws.Range("J122").Value = 10800
$ws.Range("K122").Value = 5509998
$ws.Range("L122").Value = 32400
$ws.Range("M122").Value = -5507548
$ws.Range("N122").Value = -37300
$ws.Range("H132").Value = 11463.243
$ws.Range("I132").Value = 15457.187
$ws.Range("K132").Value = 46371.561
$ws.Range("M132").Value = -43841.561

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 582205
$ws.Range("J99").Value = 4399.6523
$ws.Range("L99").Value = 4399.6523
$ws.Range("N99").Value = -7395.6523
$ws.Range("H107").Value = 1711.2
$ws.Range("I107").Value = 1743.6666
$ws.Range("J107").Value = 1662.5
$ws.Range("K107").Value = 1743.6666
$ws.Range("L107").Value = 1662.5
$ws.Range("M107").Value = 176.3334
$ws.Range("N107").Value = -5502.5

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20412616
$ws.Range("I31").Value = 76925840
$ws.Range("J31").Value = 5060.778
$ws.Range("K31").Value = 76925840
$ws.Range("L31").Value = 5060.778
$ws.Range("M31").Value = -76925545
$ws.Range("N31").Value = -5650.778
$ws.Range("H34").Value = 20412616
$ws.Range("I34").Value = 76925840
$ws.Range("J34").Value = 5060.778
$ws.Range("K34").Value = 76925840
$ws.Range("L34").Value = 5060.778
$ws.Range("M34").Value = -76925638
$ws.Range("N34").Value = -5464.778
$ws.Range("H107").Value = 536408.1
$ws.Range("I107").Value = 957812.5
$ws.Range("J107").Value = 2629.2
$ws.Range("K107").Value = 957812.5
$ws.Range("L107").Value = 2629.2
$ws.Range("M107").Value = -955892.5
$ws.Range("N107").Value = -6469.2
$ws.Range("H132").Value = 88905096
$ws.Range("I132").Value = 111114376
$ws.Range("J132").Value = 67999.664
$ws.Range("K132").Value = 333343128
$ws.Range("L132").Value = 203998.992
$ws.Range("M132").Value = -333340598
$ws.Range("N132").Value = -209058.992

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 379.25
$ws.Range("J40").Value = 1351
$ws.Range("L40").Value = 5404
$ws.Range("N40").Value = -5542
$ws.Range("H68").Value = 245001.66
$ws.Range("J68").Value = 279716.22
$ws.Range("L68").Value = 839148.6599999999
$ws.Range("N68").Value = -840770.6599999999
$ws.Range("H71").Value = 245001.66
$ws.Range("J71").Value = 279716.22
$ws.Range("L71").Value = 2517445.98
$ws.Range("N71").Value = -2525557.98
$ws.Range("H74").Value = 18875
$ws.Range("J74").Value = 21000
$ws.Range("L74").Value = 63000
$ws.Range("N74").Value = -65122
$ws.Range("H77").Value = 18875
$ws.Range("J77").Value = 21000
$ws.Range("L77").Value = 189000
$ws.Range("N77").Value = -199608
$ws.Range("H129").Value = 2508.6924
$ws.Range("I129").Value = 1197.5
$ws.Range("J129").Value = 3632.5715
$ws.Range("K129").Value = 3592.5
$ws.Range("L129").Value = 10897.7145
$ws.Range("M129").Value = 1407.5
$ws.Range("N129").Value = -20897.7145

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 50051
$ws.Range("J93").Value = 50051
$ws.Range("L93").Value = 50051
$ws.Range("N93").Value = -53795
$ws.Range("H102").Value = 10017.95
$ws.Range("I102").Value = 8797.308000000001
$ws.Range("J102").Value = 12284.857
$ws.Range("K102").Value = 8797.308000000001
$ws.Range("L102").Value = 12284.857
$ws.Range("M102").Value = -7175.308000000001
$ws.Range("N102").Value = -15528.857
$ws.Range("H122").Value = 556294
$ws.Range("I122").Value = 738765.4399999999
$ws.Range("K122").Value = 2216296.32
$ws.Range("M122").Value = -2213846.32
$ws.Range("H140").Value = 29999.924
$ws.Range("J140").Value = 29999.924
$ws.Range("L140").Value = 29999.924
$ws.Range("N140").Value = -40359.924

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3597.582
$ws.Range("I7").Value = 3017.8965
$ws.Range("K7").Value = 3017.8965
$ws.Range("M7").Value = -2905.8965
$ws.Range("H45").Value = 17520.5
$ws.Range("I45").Value = 17520.5
$ws.Range("K45").Value = 17520.5
$ws.Range("M45").Value = -17113.5
$ws.Range("H126").Value = 3597.582
$ws.Range("I126").Value = 3017.8965
$ws.Range("K126").Value = 9053.6895
$ws.Range("M126").Value = -6583.6895
$ws.Range("H131").Value = 44698
$ws.Range("J131").Value = 44698
$ws.Range("L131").Value = 44698
$ws.Range("N131").Value = -54778
$ws.Range("H139").Value = 85342
$ws.Range("J139").Value = 85342
$ws.Range("L139").Value = 85342
$ws.Range("N139").Value = -95622

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 20007694
$ws.Range("I132").Value = 8602.75
$ws.Range("K132").Value = 25808.25
$ws.Range("M132").Value = -23278.25
$ws.Range("H133").Value = 90539.664
$ws.Range("J133").Value = 90539.664
$ws.Range("L133").Value = 90539.664
$ws.Range("N133").Value = -100659.664
